$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename header columns H1 and I1
$ws.Range("I1").Value = "Air_Location"
$ws.Range("H1").Value = "Water_Temp_C"

# Update the selection to I7
$ws.Range("I7").Select()

# Set column I width to fit the new header text (mirrors Excel's bestFit column sizing)
$ws.Columns.Item(9).ColumnWidth = 10.666666666666666
